$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Rename [FINAL PHASE] keyword to [RESOLUTION PHASE]
$ws.Range("C6").Value = "[RESOLUTION PHASE]"

# Rename "specific card name" wording to "specific card tag"
$ws.Range("F10").Value = "1 <specific card tag>"
$ws.Range("F11").Value = "1 <specific card tag> and 1 <specific card tag>"
$ws.Range("F12").Value = "1 <specific card tag> or 1 <specific card tag>"

# Update the "on attack" keyword to "on ambush"
$ws.Range("B3").Value = "(on ambush)"

# Column C needs to widen (bestFit) now that its longest entry changed
$ws.Columns.Item(3).AutoFit() | Out-Null

# Move the active selection to B3
$ws.Activate()
$ws.Range("B3").Select()
